# Add "CIRP Activity" entry to the metadata mapping (Zuordnung) sheet and
# to the activity list on the Listen sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Zuordnung")
$ws2 = $wb.Worksheets.Item("Listen")

# --- Zuordnung!A9:O9 : new "CIRP Activity" row -----------------------------
$ws1.Range("A9").Value = "CIRP Activity"
$ws1.Range("B9").Value = "Activity"
$ws1.Range("C9").Value = "activity"
$ws1.Range("D9").Value = "dropdown"
$ws1.Range("F9").Value = "default"
$ws1.Range("G9:O9").Value = "Required"

# Extend the conditional formatting that marks "dropdown"/"suggestions"
# definitions so it also covers the new row's Required-cells (H9:J9).
$rngCf = $ws1.Range("D2:F1048576,H9:J9")
$cfRules = $ws1.Range("D2").FormatConditions
$cfRules.Item(1).ModifyAppliesToRange($rngCf) | Out-Null
$cfRules.Item(2).ModifyAppliesToRange($rngCf) | Out-Null

# --- Listen!G3:G4 : new "activity" lookup column ---------------------------
$ws2.Range("G3").Value = "activity"
$ws2.Range("G4").Value = "Tool Wear Monitoring"
$ws2.Columns.Item(7).ColumnWidth = 18.14

# --- restore the selections / active sheet as saved in the workbook -------
$ws1.Activate() | Out-Null
$ws1.Range("C31").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("G11").Select() | Out-Null
